$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value  = 153
$ws.Range("I2").Value  = 343
$ws.Range("J2").Value  = 1438
$ws.Range("K2").Value  = 2
$ws.Range("L2").Value  = 395
$ws.Range("M2").Value  = 25
$ws.Range("N2").Value  = 220
$ws.Range("O2").Value  = 2
$ws.Range("R2").Value  = 17
$ws.Range("S2").Value  = 173
$ws.Range("T2").Value  = 229
$ws.Range("U2").Value  = 23
$ws.Range("V2").Value  = 2193
$ws.Range("W2").Value  = 1
$ws.Range("X2").Value  = 2184
$ws.Range("Y2").Value  = 3
$ws.Range("Z2").Value  = 23
$ws.Range("AA2").Value = 21
